# Add a new task row "Mob Special Attack/Defense, Magic Resistance"
# This inserts a new row at row 14 (pushing all following rows down by one),
# fills in the new task name in column C, and restores the cell selection
# to C13 (the cell above the newly inserted row), matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row above the current row 14 (shifts rows 14+ down by 1)
$ws.Rows.Item(14).Insert()

# Populate the new row's Task column with the new task name
$ws.Range("C14").Value = "Mob Special Attack/Defense, Magic Resistance"

# Restore the selected cell to C13
$ws.Range("C13").Select()
